$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 16.38931533333333
$ws.Cells.Item(2, 8).Value = 49.167946
$ws.Cells.Item(2, 9).Value = 0.1272611691689642
$ws.Cells.Item(2, 10).Value = 0.1272611691689643
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 24.455837
$ws.Cells.Item(2, 14).Value = 73.36751100000001
$ws.Cells.Item(2, 15).Value = 0.1553502885444182
$ws.Cells.Item(2, 16).Value = 0.1553502885444182
$ws.Cells.Item(2, 17).Value = 400.8144243336007
$ws.Cells.Item(2, 18).Value = 3607.329819002407
$ws.Cells.Item(2, 19).Value = 0.01977005935089861
$ws.Cells.Item(2, 20).Value = 0.01977005935089862

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 16.38931533333333
$ws.Cells.Item(3, 8).Value = 49.167946
$ws.Cells.Item(3, 9).Value = 0.1272611691689642
$ws.Cells.Item(3, 10).Value = 0.1272611691689643
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 33.36516466666667
$ws.Cells.Item(3, 14).Value = 100.095494
$ws.Cells.Item(3, 15).Value = 0.2119448194841458
$ws.Cells.Item(3, 16).Value = 0.2119448194841458
$ws.Cells.Item(3, 17).Value = 546.8322048705915
$ws.Cells.Item(3, 18).Value = 4921.489843835324
$ws.Cells.Item(3, 19).Value = 0.02697234552685747
$ws.Cells.Item(3, 20).Value = 0.02697234552685747

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 16.38931533333333
$ws.Cells.Item(4, 8).Value = 49.167946
$ws.Cells.Item(4, 9).Value = 0.1272611691689642
$ws.Cells.Item(4, 10).Value = 0.1272611691689643
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 43.331795
$ws.Cells.Item(4, 14).Value = 129.995385
$ws.Cells.Item(4, 15).Value = 0.2752556314632608
$ws.Cells.Item(4, 16).Value = 0.2752556314632608
$ws.Cells.Item(4, 17).Value = 710.1784522143566
$ws.Cells.Item(4, 18).Value = 6391.60606992921
$ws.Cells.Item(4, 19).Value = 0.03502935348035611
$ws.Cells.Item(4, 20).Value = 0.03502935348035611

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 16.38931533333333
$ws.Cells.Item(5, 8).Value = 49.167946
$ws.Cells.Item(5, 9).Value = 0.1272611691689642
$ws.Cells.Item(5, 10).Value = 0.1272611691689643
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 45.91844699999999
$ws.Cells.Item(5, 14).Value = 137.755341
$ws.Cells.Item(5, 15).Value = 0.2916867654524183
$ws.Cells.Item(5, 16).Value = 0.2916867654524183
$ws.Cells.Item(5, 17).Value = 752.5719074999538
$ws.Cells.Item(5, 18).Value = 6773.147167499586
$ws.Cells.Item(5, 19).Value = 0.0371203988025882
$ws.Cells.Item(5, 20).Value = 0.0371203988025882

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 16.38931533333333
$ws.Cells.Item(6, 8).Value = 49.167946
$ws.Cells.Item(6, 9).Value = 0.1272611691689642
$ws.Cells.Item(6, 10).Value = 0.1272611691689643
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 10.35258366666667
$ws.Cells.Item(6, 14).Value = 31.057751
$ws.Cells.Item(6, 15).Value = 0.06576249505575693
$ws.Cells.Item(6, 16).Value = 0.06576249505575693
$ws.Cells.Item(6, 17).Value = 169.6717582277162
$ws.Cells.Item(6, 18).Value = 1527.045824049446
$ws.Cells.Item(6, 19).Value = 0.008369012008263857
$ws.Cells.Item(6, 20).Value = 0.008369012008263859

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 23.071008
$ws.Cells.Item(7, 8).Value = 69.213024
$ws.Cells.Item(7, 9).Value = 0.1791437526383466
$ws.Cells.Item(7, 10).Value = 0.1791437526383466
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 24.455837
$ws.Cells.Item(7, 14).Value = 73.36751100000001
$ws.Cells.Item(7, 15).Value = 0.1553502885444182
$ws.Cells.Item(7, 16).Value = 0.1553502885444182
$ws.Cells.Item(7, 17).Value = 564.2208110736962
$ws.Cells.Item(7, 18).Value = 5077.987299663265
$ws.Cells.Item(7, 19).Value = 0.02783003366329703
$ws.Cells.Item(7, 20).Value = 0.02783003366329703

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 23.071008
$ws.Cells.Item(8, 8).Value = 69.213024
$ws.Cells.Item(8, 9).Value = 0.1791437526383466
$ws.Cells.Item(8, 10).Value = 0.1791437526383466
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 33.36516466666667
$ws.Cells.Item(8, 14).Value = 100.095494
$ws.Cells.Item(8, 15).Value = 0.2119448194841458
$ws.Cells.Item(8, 16).Value = 0.2119448194841458
$ws.Cells.Item(8, 17).Value = 769.767980945984
$ws.Cells.Item(8, 18).Value = 6927.911828513857
$ws.Cells.Item(8, 19).Value = 0.03796859031464684
$ws.Cells.Item(8, 20).Value = 0.03796859031464684

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 23.071008
$ws.Cells.Item(9, 8).Value = 69.213024
$ws.Cells.Item(9, 9).Value = 0.1791437526383466
$ws.Cells.Item(9, 10).Value = 0.1791437526383466
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 43.331795
$ws.Cells.Item(9, 14).Value = 129.995385
$ws.Cells.Item(9, 15).Value = 0.2752556314632608
$ws.Cells.Item(9, 16).Value = 0.2752556314632608
$ws.Cells.Item(9, 17).Value = 999.7081890993601
$ws.Cells.Item(9, 18).Value = 8997.373701894241
$ws.Cells.Item(9, 19).Value = 0.04931032675516629
$ws.Cells.Item(9, 20).Value = 0.04931032675516629

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 23.071008
$ws.Cells.Item(10, 8).Value = 69.213024
$ws.Cells.Item(10, 9).Value = 0.1791437526383466
$ws.Cells.Item(10, 10).Value = 0.1791437526383466
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 45.91844699999999
$ws.Cells.Item(10, 14).Value = 137.755341
$ws.Cells.Item(10, 15).Value = 0.2916867654524183
$ws.Cells.Item(10, 16).Value = 0.2916867654524183
$ws.Cells.Item(10, 17).Value = 1059.384858084576
$ws.Cells.Item(10, 18).Value = 9534.463722761184
$ws.Cells.Item(10, 19).Value = 0.05225386175808745
$ws.Cells.Item(10, 20).Value = 0.05225386175808745

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 23.071008
$ws.Cells.Item(11, 8).Value = 69.213024
$ws.Cells.Item(11, 9).Value = 0.1791437526383466
$ws.Cells.Item(11, 10).Value = 0.1791437526383466
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 10.35258366666667
$ws.Cells.Item(11, 14).Value = 31.057751
$ws.Cells.Item(11, 15).Value = 0.06576249505575693
$ws.Cells.Item(11, 16).Value = 0.06576249505575693
$ws.Cells.Item(11, 17).Value = 238.844540594336
$ws.Cells.Item(11, 18).Value = 2149.600865349024
$ws.Cells.Item(11, 19).Value = 0.01178094014714901
$ws.Cells.Item(11, 20).Value = 0.01178094014714901

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 37.292974
$ws.Cells.Item(12, 8).Value = 111.878922
$ws.Cells.Item(12, 9).Value = 0.2895757007844777
$ws.Cells.Item(12, 10).Value = 0.2895757007844777
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 24.455837
$ws.Cells.Item(12, 14).Value = 73.36751100000001
$ws.Cells.Item(12, 15).Value = 0.1553502885444182
$ws.Cells.Item(12, 16).Value = 0.1553502885444182
$ws.Cells.Item(12, 17).Value = 912.0308933892381
$ws.Cells.Item(12, 18).Value = 8208.278040503143
$ws.Cells.Item(12, 19).Value = 0.04498566867232072
$ws.Cells.Item(12, 20).Value = 0.04498566867232072

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 37.292974
$ws.Cells.Item(13, 8).Value = 111.878922
$ws.Cells.Item(13, 9).Value = 0.2895757007844777
$ws.Cells.Item(13, 10).Value = 0.2895757007844777
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 33.36516466666667
$ws.Cells.Item(13, 14).Value = 100.095494
$ws.Cells.Item(13, 15).Value = 0.2119448194841458
$ws.Cells.Item(13, 16).Value = 0.2119448194841458
$ws.Cells.Item(13, 17).Value = 1244.286218419719
$ws.Cells.Item(13, 18).Value = 11198.57596577747
$ws.Cells.Item(13, 19).Value = 0.06137406962976114
$ws.Cells.Item(13, 20).Value = 0.06137406962976114

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 37.292974
$ws.Cells.Item(14, 8).Value = 111.878922
$ws.Cells.Item(14, 9).Value = 0.2895757007844777
$ws.Cells.Item(14, 10).Value = 0.2895757007844777
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 43.331795
$ws.Cells.Item(14, 14).Value = 129.995385
$ws.Cells.Item(14, 15).Value = 0.2752556314632608
$ws.Cells.Item(14, 16).Value = 0.2752556314632608
$ws.Cells.Item(14, 17).Value = 1615.97150430833
$ws.Cells.Item(14, 18).Value = 14543.74353877497
$ws.Cells.Item(14, 19).Value = 0.07970734237584767
$ws.Cells.Item(14, 20).Value = 0.07970734237584766

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 37.292974
$ws.Cells.Item(15, 8).Value = 111.878922
$ws.Cells.Item(15, 9).Value = 0.2895757007844777
$ws.Cells.Item(15, 10).Value = 0.2895757007844777
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 45.91844699999999
$ws.Cells.Item(15, 14).Value = 137.755341
$ws.Cells.Item(15, 15).Value = 0.2916867654524183
$ws.Cells.Item(15, 16).Value = 0.2916867654524183
$ws.Cells.Item(15, 17).Value = 1712.435450091378
$ws.Cells.Item(15, 18).Value = 15411.9190508224
$ws.Cells.Item(15, 19).Value = 0.0844653995154416
$ws.Cells.Item(15, 20).Value = 0.0844653995154416

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 37.292974
$ws.Cells.Item(16, 8).Value = 111.878922
$ws.Cells.Item(16, 9).Value = 0.2895757007844777
$ws.Cells.Item(16, 10).Value = 0.2895757007844777
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 10.35258366666667
$ws.Cells.Item(16, 14).Value = 31.057751
$ws.Cells.Item(16, 15).Value = 0.06576249505575693
$ws.Cells.Item(16, 16).Value = 0.06576249505575693
$ws.Cells.Item(16, 17).Value = 386.0786335138246
$ws.Cells.Item(16, 18).Value = 3474.707701624422
$ws.Cells.Item(16, 19).Value = 0.01904322059110656
$ws.Cells.Item(16, 20).Value = 0.01904322059110656

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 43.83143633333333
$ws.Cells.Item(17, 8).Value = 131.494309
$ws.Cells.Item(17, 9).Value = 0.3403461170089362
$ws.Cells.Item(17, 10).Value = 0.3403461170089362
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 24.455837
$ws.Cells.Item(17, 14).Value = 73.36751100000001
$ws.Cells.Item(17, 15).Value = 0.1553502885444182
$ws.Cells.Item(17, 16).Value = 0.1553502885444182
$ws.Cells.Item(17, 17).Value = 1071.934462443878
$ws.Cells.Item(17, 18).Value = 9647.4101619949
$ws.Cells.Item(17, 19).Value = 0.05287286748231056
$ws.Cells.Item(17, 20).Value = 0.05287286748231056

$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 7).Value = 43.83143633333333
$ws.Cells.Item(18, 8).Value = 131.494309
$ws.Cells.Item(18, 9).Value = 0.3403461170089362
$ws.Cells.Item(18, 10).Value = 0.3403461170089362
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 13).Value = 33.36516466666667
$ws.Cells.Item(18, 14).Value = 100.095494
$ws.Cells.Item(18, 15).Value = 0.2119448194841458
$ws.Cells.Item(18, 16).Value = 0.2119448194841458
$ws.Cells.Item(18, 17).Value = 1462.443090838183
$ws.Cells.Item(18, 18).Value = 13161.98781754364
$ws.Cells.Item(18, 19).Value = 0.07213459633158895
$ws.Cells.Item(18, 20).Value = 0.07213459633158895

$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 7).Value = 43.83143633333333
$ws.Cells.Item(19, 8).Value = 131.494309
$ws.Cells.Item(19, 9).Value = 0.3403461170089362
$ws.Cells.Item(19, 10).Value = 0.3403461170089362
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 13).Value = 43.331795
$ws.Cells.Item(19, 14).Value = 129.995385
$ws.Cells.Item(19, 15).Value = 0.2752556314632608
$ws.Cells.Item(19, 16).Value = 0.2752556314632608
$ws.Cells.Item(19, 17).Value = 1899.294813751551
$ws.Cells.Item(19, 18).Value = 17093.65332376396
$ws.Cells.Item(19, 19).Value = 0.09368218535336359
$ws.Cells.Item(19, 20).Value = 0.09368218535336356

$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 7).Value = 43.83143633333333
$ws.Cells.Item(20, 8).Value = 131.494309
$ws.Cells.Item(20, 9).Value = 0.3403461170089362
$ws.Cells.Item(20, 10).Value = 0.3403461170089362
$ws.Cells.Item(20, 11).Value = 3
$ws.Cells.Item(20, 13).Value = 45.91844699999999
$ws.Cells.Item(20, 14).Value = 137.755341
$ws.Cells.Item(20, 15).Value = 0.2916867654524183
$ws.Cells.Item(20, 16).Value = 0.2916867654524183
$ws.Cells.Item(20, 17).Value = 2012.671486206041
$ws.Cells.Item(20, 18).Value = 18114.04337585437
$ws.Cells.Item(20, 19).Value = 0.09927445800462688
$ws.Cells.Item(20, 20).Value = 0.09927445800462688

$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 7).Value = 43.83143633333333
$ws.Cells.Item(21, 8).Value = 131.494309
$ws.Cells.Item(21, 9).Value = 0.3403461170089362
$ws.Cells.Item(21, 10).Value = 0.3403461170089362
$ws.Cells.Item(21, 11).Value = 3
$ws.Cells.Item(21, 13).Value = 10.35258366666667
$ws.Cells.Item(21, 14).Value = 31.057751
$ws.Cells.Item(21, 15).Value = 0.06576249505575693
$ws.Cells.Item(21, 16).Value = 0.06576249505575693
$ws.Cells.Item(21, 17).Value = 453.7686118710065
$ws.Cells.Item(21, 18).Value = 4083.917506839059
$ws.Cells.Item(21, 19).Value = 0.02238200983704624
$ws.Cells.Item(21, 20).Value = 0.02238200983704624

$ws.Cells.Item(22, 5).Value = 3
$ws.Cells.Item(22, 7).Value = 8.200153666666667
$ws.Cells.Item(22, 8).Value = 24.600461
$ws.Cells.Item(22, 9).Value = 0.06367326039927532
$ws.Cells.Item(22, 10).Value = 0.06367326039927532
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 13).Value = 24.455837
$ws.Cells.Item(22, 14).Value = 73.36751100000001
$ws.Cells.Item(22, 15).Value = 0.1553502885444182
$ws.Cells.Item(22, 16).Value = 0.1553502885444182
$ws.Cells.Item(22, 17).Value = 200.5416214469523
$ws.Cells.Item(22, 18).Value = 1804.874593022571
$ws.Cells.Item(22, 19).Value = 0.009891659375591298
$ws.Cells.Item(22, 20).Value = 0.009891659375591298

$ws.Cells.Item(23, 5).Value = 3
$ws.Cells.Item(23, 7).Value = 8.200153666666667
$ws.Cells.Item(23, 8).Value = 24.600461
$ws.Cells.Item(23, 9).Value = 0.06367326039927532
$ws.Cells.Item(23, 10).Value = 0.06367326039927532
$ws.Cells.Item(23, 11).Value = 3
$ws.Cells.Item(23, 13).Value = 33.36516466666667
$ws.Cells.Item(23, 14).Value = 100.095494
$ws.Cells.Item(23, 15).Value = 0.2119448194841458
$ws.Cells.Item(23, 16).Value = 0.2119448194841458
$ws.Cells.Item(23, 17).Value = 273.5994773803038
$ws.Cells.Item(23, 18).Value = 2462.395296422734
$ws.Cells.Item(23, 19).Value = 0.01349521768129142
$ws.Cells.Item(23, 20).Value = 0.01349521768129142

$ws.Cells.Item(24, 5).Value = 3
$ws.Cells.Item(24, 7).Value = 8.200153666666667
$ws.Cells.Item(24, 8).Value = 24.600461
$ws.Cells.Item(24, 9).Value = 0.06367326039927532
$ws.Cells.Item(24, 10).Value = 0.06367326039927532
$ws.Cells.Item(24, 11).Value = 3
$ws.Cells.Item(24, 13).Value = 43.331795
$ws.Cells.Item(24, 14).Value = 129.995385
$ws.Cells.Item(24, 15).Value = 0.2752556314632608
$ws.Cells.Item(24, 16).Value = 0.2752556314632608
$ws.Cells.Item(24, 17).Value = 355.3273776524983
$ws.Cells.Item(24, 18).Value = 3197.946398872485
$ws.Cells.Item(24, 19).Value = 0.01752642349852717
$ws.Cells.Item(24, 20).Value = 0.01752642349852716

$ws.Cells.Item(25, 5).Value = 3
$ws.Cells.Item(25, 7).Value = 8.200153666666667
$ws.Cells.Item(25, 8).Value = 24.600461
$ws.Cells.Item(25, 9).Value = 0.06367326039927532
$ws.Cells.Item(25, 10).Value = 0.06367326039927532
$ws.Cells.Item(25, 11).Value = 3
$ws.Cells.Item(25, 13).Value = 45.91844699999999
$ws.Cells.Item(25, 14).Value = 137.755341
$ws.Cells.Item(25, 15).Value = 0.2916867654524183
$ws.Cells.Item(25, 16).Value = 0.2916867654524183
$ws.Cells.Item(25, 17).Value = 376.538321534689
$ws.Cells.Item(25, 18).Value = 3388.844893812201
$ws.Cells.Item(25, 19).Value = 0.01857264737167417
$ws.Cells.Item(25, 20).Value = 0.01857264737167417

$ws.Cells.Item(26, 5).Value = 3
$ws.Cells.Item(26, 7).Value = 8.200153666666667
$ws.Cells.Item(26, 8).Value = 24.600461
$ws.Cells.Item(26, 9).Value = 0.06367326039927532
$ws.Cells.Item(26, 10).Value = 0.06367326039927532
$ws.Cells.Item(26, 11).Value = 3
$ws.Cells.Item(26, 13).Value = 10.35258366666667
$ws.Cells.Item(26, 14).Value = 31.057751
$ws.Cells.Item(26, 15).Value = 0.06576249505575693
$ws.Cells.Item(26, 16).Value = 0.06576249505575693
$ws.Cells.Item(26, 17).Value = 84.89277691369011
$ws.Cells.Item(26, 18).Value = 764.034992223211
$ws.Cells.Item(26, 19).Value = 0.004187312472191266
$ws.Cells.Item(26, 20).Value = 0.004187312472191266
